# Daily attendance processing - 2026-02-18 16:43:46 UTC
# Reorders the comma-separated "Recorded By" values in column G so that the
# "2025/2026" academic-year token is listed after the other recorder(s)
# instead of before them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = "2022/2023, 2025/2026"
    22 = "2024/2025, 2025/2026"
    23 = "2023/2024, 2022/2023, 2025/2026"
    24 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    27 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    28 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    31 = "2022/2023, 2025/2026"
    50 = "2024/2025, 2025/2026"
    51 = "2023/2024, 2022/2023, 2025/2026"
    52 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    55 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
    56 = "neveen.nashaat@med.asu.edu.eg, 2025/2026"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
